# Swap the presentation's applied theme color palette from the custom
# "Integral" scheme back to the stock Office "Office Theme" scheme.
#
# The deck currently carries two theme parts: the one driving the slide
# master/presentation (Integral: greens/teals) and a second one used only
# by the notes master (the stock Office Theme: blues/oranges). This change
# swaps which palette is which by rewriting the 12 theme colors exposed
# through the Slide Master's ThemeColorScheme to the stock Office values.
# (OLE RGB is packed as 0x00BBGGRR, i.e. R + G*256 + B*65536.)

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Index order (matches document order dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) -> target "Office Theme" RGB values.
$cs.Item(1).RGB  = 0            # dk1      000000
$cs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388      # dk2      44546A
$cs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501      # accent2  ED7D31
$cs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$cs.Item(8).RGB  = 49407        # accent4  FFC000
$cs.Item(9).RGB  = 12874308     # accent5  4472C4
$cs.Item(10).RGB = 4697456      # accent6  70AD47
$cs.Item(11).RGB = 12673797     # hlink    0563C1
$cs.Item(12).RGB = 7491477      # folHlink 954F72
